# Commit: "Tried to implement Penality Reward System (unfinished)"
# Net effect observed in the diff: a data row was removed from each sheet
# (the remaining rows shift up by one), reducing each sheet's used range
# by one row.
$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: remove row 10 (45102.99999999999 / 40),
# shifting the former rows 11-12 up to become rows 10-11.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(10).Delete()

# "Monthly Trend" sheet: remove row 5 (45107.99999999999 / 40),
# shifting the former row 6 up to become row 5.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows.Item(5).Delete()
